$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.550.25'
$ws.Cells.Item(2, 5).Value = '  -2.29%  '
$ws.Cells.Item(3, 4).Value = '2.562.45'
$ws.Cells.Item(3, 5).Value = '  -3.78%  '
$ws.Cells.Item(4, 5).Value = '  +0.00%  '
$ws.Cells.Item(5, 4).Value = '''520.59'
$ws.Cells.Item(5, 5).Value = '  -0.77%  '
$ws.Cells.Item(6, 4).Value = '''143.71'
$ws.Cells.Item(6, 5).Value = '  -0.37%  '
$ws.Cells.Item(7, 5).Value = '  -0.11%  '
$ws.Cells.Item(8, 4).Value = '''0.562'
$ws.Cells.Item(8, 5).Value = '  -1.50%  '
$ws.Cells.Item(9, 4).Value = '2.573.54'
$ws.Cells.Item(10, 4).Value = '''6.64'
$ws.Cells.Item(10, 5).Value = '  -4.03%  '
$ws.Cells.Item(11, 5).Value = '  -2.37%  '
$ws.Cells.Item(12, 5).Value = '  -2.91%  '
$ws.Cells.Item(13, 5).Value = '  -0.33%  '
$ws.Cells.Item(14, 4).Value = '3.010.83'
$ws.Cells.Item(14, 5).Value = '  -3.92%  '
$ws.Cells.Item(15, 4).Value = '57.512.91'
$ws.Cells.Item(15, 5).Value = '  -2.31%  '
$ws.Cells.Item(16, 4).Value = '''20.16'
$ws.Cells.Item(16, 5).Value = '  -4.16%  '
$ws.Cells.Item(17, 5).Value = '  -2.61%  '
$ws.Cells.Item(18, 4).Value = '2.569.51'
$ws.Cells.Item(18, 5).Value = '  -3.64%  '
$ws.Cells.Item(19, 4).Value = '''335.18'
$ws.Cells.Item(19, 5).Value = '  -1.27%  '
$ws.Cells.Item(20, 5).Value = '  -2.26%  '
$ws.Cells.Item(21, 4).Value = '''10.18'
$ws.Cells.Item(21, 5).Value = '  -2.32%  '
$ws.Cells.Item(22, 4).Value = '''6.22'
$ws.Cells.Item(22, 5).Value = '  -2.49%  '
$ws.Cells.Item(23, 5).Value = '  -0.10%  '
$ws.Cells.Item(24, 4).Value = '''65.19'
$ws.Cells.Item(24, 5).Value = '  +1.59%  '
$ws.Cells.Item(25, 5).Value = '  -0.74%  '
$ws.Cells.Item(26, 4).Value = '''0.402'
$ws.Cells.Item(26, 5).Value = '  -5.08%  '
$ws.Cells.Item(27, 4).Value = '''0.997'
$ws.Cells.Item(27, 5).Value = '  -0.14%  '
$ws.Cells.Item(28, 4).Value = '2.678.98'
$ws.Cells.Item(28, 5).Value = '  -3.71%  '
$ws.Cells.Item(29, 5).Value = '  -2.86%  '
$ws.Cells.Item(30, 4).Value = '0.0₃0748'
$ws.Cells.Item(30, 5).Value = '  -7.28%  '
$ws.Cells.Item(31, 5).Value = '  +0.00%  '
$ws.Cells.Item(32, 5).Value = '  -7.04%  '
$ws.Cells.Item(33, 5).Value = '  -0.92%  '
$ws.Cells.Item(34, 4).Value = '''18.63'
$ws.Cells.Item(34, 5).Value = '  -1.47%  '
$ws.Cells.Item(35, 4).Value = '''148.73'
$ws.Cells.Item(35, 5).Value = '  -1.48%  '
$ws.Cells.Item(36, 4).Value = '''4.03'
$ws.Cells.Item(36, 5).Value = '  -3.15%  '
$ws.Cells.Item(37, 5).Value = '  -4.18%  '
$ws.Cells.Item(38, 4).Value = '''0.843'
$ws.Cells.Item(38, 5).Value = '  -9.32%  '
$ws.Cells.Item(39, 4).Value = '''36.07'
$ws.Cells.Item(39, 5).Value = '  -1.96%  '
$ws.Cells.Item(40, 4).Value = '''0.831'
$ws.Cells.Item(40, 5).Value = '  -5.10%  '
$ws.Cells.Item(41, 4).Value = '''1.43'
$ws.Cells.Item(41, 5).Value = '  -1.40%  '
$ws.Cells.Item(42, 5).Value = '  -2.46%  '
$ws.Cells.Item(43, 5).Value = '  -0.17%  '
$ws.Cells.Item(44, 4).Value = '''268.11'
$ws.Cells.Item(44, 5).Value = '  -2.90%  '
$ws.Cells.Item(45, 5).Value = '  -1.34%  '
$ws.Cells.Item(47, 4).Value = '''0.587'
$ws.Cells.Item(47, 5).Value = '  -4.24%  '
$ws.Cells.Item(48, 4).Value = '''18.79'
$ws.Cells.Item(48, 5).Value = '  -4.70%  '
$ws.Cells.Item(49, 5).Value = '  -2.54%  '
$ws.Cells.Item(50, 4).Value = '1.969.64'
$ws.Cells.Item(51, 4).Value = '''4.63'
$ws.Cells.Item(51, 5).Value = '  -2.30%  '
